# Extend the "Payment process..." closing paragraph with the extra
# sentence about expections-first testing, then append a blank paragraph
# and a new closing paragraph describing the testing approach.

$d = $word.ActiveDocument

# Locate the final paragraph ("Payment process is the next one to be
# developed") and append the extra clause to its existing run/text,
# collapsing the range to its end first so we only insert, not replace,
# keeping the run's (empty) rPr intact.
$paymentPara = $d.Paragraphs.Last
$paymentRange = $paymentPara.Range
$paymentRange.Collapse(0)
$paymentRange.InsertAfter(" – the first tests for these are the expections as they are the simpler option to handle ")

# Add a new, empty paragraph after it (mirrors the blank-line spacing
# used elsewhere in the document).
$blankPara = $d.Paragraphs.Last
$blankPara.Range.InsertParagraphAfter()

# Add the final paragraph describing the testing approach.
$newPara = $d.Paragraphs.Last
$newPara.Range.InsertParagraphAfter()
$testingPara = $d.Paragraphs.Last
$testingPara.Range.Text = "Testing approach – start from edge cases then go into generic tests – this ensures better code coverage and that the system is tested fully"
